# Generate Report for Archive
# The localization status report was regenerated; the status of the file
# 6603f24f-90f7-44ae-8799-dd3d11673c6d.md moved on from "Ready for handoff"
# to "In Translation" on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B4").Value = "In Translation"
$wsOverview.Range("C4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C4").Value = "In Translation"
